$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.010513544082642
$ws.Range("B1").Value = 2.122204303741455
$ws.Range("C1").Value = 6.17284631729126
$ws.Range("D1").Value = 1.381744861602783
$ws.Range("E1").Value = 1.313313961029053
